$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2879077579581804
$ws.Cells.Item(2, 3).Value = 0.06802732989288529
$ws.Cells.Item(2, 4).Value = 0.03218597935602219
$ws.Cells.Item(2, 6).Value = 0.4429201865440646
$ws.Cells.Item(2, 7).Value = 0.2890509517787407
$ws.Cells.Item(2, 8).Value = 0.475572796183819
$ws.Cells.Item(2, 11).Value = 0.2798556452117396
$ws.Cells.Item(2, 14).Value = 1.09904131137565
$ws.Cells.Item(2, 15).Value = 1.454800420745968
$ws.Cells.Item(3, 2).Value = 0.2527692589264632
$ws.Cells.Item(3, 3).Value = 0.06480976267991423
$ws.Cells.Item(3, 4).Value = 0.02846685834177265
$ws.Cells.Item(3, 6).Value = 0.4444979314815249
$ws.Cells.Item(3, 7).Value = 0.2916835176808945
$ws.Cells.Item(3, 8).Value = 0.4797159895050598
$ws.Cells.Item(3, 11).Value = 0.244258678872626
$ws.Cells.Item(3, 14).Value = 1.10666087772519
$ws.Cells.Item(3, 15).Value = 1.468864360362559
$ws.Cells.Item(4, 2).Value = 0.2311499199164757
$ws.Cells.Item(4, 3).Value = 0.06282655412506699
$ws.Cells.Item(4, 4).Value = 0.02617026664467659
$ws.Cells.Item(4, 6).Value = 0.4457613045544093
$ws.Cells.Item(4, 7).Value = 0.293525649720749
$ws.Cells.Item(4, 8).Value = 0.482461176507833
$ws.Cells.Item(4, 11).Value = 0.222311113264027
$ws.Cells.Item(4, 14).Value = 1.111719578907284
$ws.Cells.Item(4, 15).Value = 1.478392936267312
$ws.Cells.Item(5, 2).Value = 0.2223294096006896
$ws.Cells.Item(5, 3).Value = 0.06201654191296058
$ws.Cells.Item(5, 4).Value = 0.02523116518965196
$ws.Cells.Item(5, 6).Value = 0.4463502547271432
$ws.Cells.Item(5, 7).Value = 0.2943330491058944
$ws.Cells.Item(5, 8).Value = 0.4836305071698206
$ws.Cells.Item(5, 11).Value = 0.2133451020998081
$ws.Cells.Item(5, 14).Value = 1.113876774381353
$ws.Cells.Item(5, 15).Value = 1.482500507363326
$ws.Cells.Item(6, 2).Value = 0.2208641582942334
$ws.Cells.Item(6, 3).Value = 0.06188193115460194
$ws.Cells.Item(6, 4).Value = 0.02507503518282306
$ws.Cells.Item(6, 6).Value = 0.4464525266552002
$ws.Cells.Item(6, 7).Value = 0.2944705413154693
$ws.Cells.Item(6, 8).Value = 0.4838277335906085
$ws.Cells.Item(6, 11).Value = 0.2118549811828672
$ws.Cells.Item(6, 14).Value = 1.114240760539623
$ws.Cells.Item(6, 15).Value = 1.483196131093507
$ws.Cells.Item(7, 2).Value = 0.2310310048970052
$ws.Cells.Item(7, 3).Value = 0.0628156373620854
$ws.Cells.Item(7, 4).Value = 0.02615761456365817
$ws.Cells.Item(7, 6).Value = 0.4457689472146242
$ws.Cells.Item(7, 7).Value = 0.2935363090113796
$ws.Cells.Item(7, 8).Value = 0.4824767413996511
$ws.Cells.Item(7, 11).Value = 0.2221902835229059
$ws.Cells.Item(7, 14).Value = 1.111748283823232
$ws.Cells.Item(7, 15).Value = 1.478447423037821
$ws.Cells.Item(8, 2).Value = 0.2758015479604694
$ws.Cells.Item(8, 3).Value = 0.06691953589825061
$ws.Cells.Item(8, 4).Value = 0.03090636495293353
$ws.Cells.Item(8, 6).Value = 0.4434030582373651
$ws.Cells.Item(8, 7).Value = 0.2899117839200898
$ws.Cells.Item(8, 8).Value = 0.4769596233576152
$ws.Cells.Item(8, 11).Value = 0.26760105588653
$ws.Cells.Item(8, 14).Value = 1.101589684275503
$ws.Cells.Item(8, 15).Value = 1.459464259881273
$ws.Cells.Item(9, 2).Value = 0.3632215908014302
$ws.Cells.Item(9, 3).Value = 0.07490397385519998
$ws.Cells.Item(9, 4).Value = 0.04011313934142891
$ws.Cells.Item(9, 6).Value = 0.4411005722768522
$ws.Cells.Item(9, 7).Value = 0.2845973197197651
$ws.Cells.Item(9, 8).Value = 0.4677356622538227
$ws.Cells.Item(9, 11).Value = 0.3559056305408888
$ws.Cells.Item(9, 14).Value = 1.084680870394237
$ws.Cells.Item(9, 15).Value = 1.429327152602568
$ws.Cells.Item(10, 2).Value = 0.4271936722313683
$ws.Cells.Item(10, 3).Value = 0.08072817216086037
$ws.Cells.Item(10, 4).Value = 0.04681084475470243
$ws.Cells.Item(10, 6).Value = 0.4408333360433403
$ws.Cells.Item(10, 7).Value = 0.2817892223437681
$ws.Cells.Item(10, 8).Value = 0.4619288547312621
$ws.Cells.Item(10, 11).Value = 0.420301089892348
$ws.Cells.Item(10, 14).Value = 1.07408762714482
$ws.Cells.Item(10, 15).Value = 1.41150914421705
$ws.Cells.Item(11, 2).Value = 0.4562354036945351
$ws.Cells.Item(11, 3).Value = 0.08336795348238013
$ws.Cells.Item(11, 4).Value = 0.04984292865127316
$ws.Cells.Item(11, 6).Value = 0.4410210296917541
$ws.Cells.Item(11, 7).Value = 0.2807505314800594
$ws.Cells.Item(11, 8).Value = 0.4594973230341282
$ws.Cells.Item(11, 11).Value = 0.4494861744080367
$ws.Cells.Item(11, 14).Value = 1.069664420564919
$ws.Cells.Item(11, 15).Value = 1.404342633757466
$ws.Cells.Item(12, 2).Value = 0.4672236074483465
$ws.Cells.Item(12, 3).Value = 0.08436610917541998
$ws.Cells.Item(12, 4).Value = 0.05098892985689929
$ws.Cells.Item(12, 6).Value = 0.4411365588302161
$ws.Cells.Item(12, 7).Value = 0.2803915878054113
$ws.Cells.Item(12, 8).Value = 0.4586067282584381
$ws.Cells.Item(12, 11).Value = 0.4605215767024049
$ws.Cells.Item(12, 14).Value = 1.068046268141224
$ws.Cells.Item(12, 15).Value = 1.401763928135765
$ws.Cells.Item(13, 2).Value = 0.4648575248252769
$ws.Cells.Item(13, 3).Value = 0.08415120524902875
$ws.Cells.Item(13, 4).Value = 0.05074221611329222
$ws.Cells.Item(13, 6).Value = 0.441109700770852
$ws.Cells.Item(13, 7).Value = 0.2804673625175838
$ws.Cells.Item(13, 8).Value = 0.4587971922215388
$ws.Cells.Item(13, 11).Value = 0.4581456442626859
$ws.Cells.Item(13, 14).Value = 1.068392240441071
$ws.Cells.Item(13, 15).Value = 1.402313288982072
$ws.Cells.Item(14, 2).Value = 0.4571395997047318
$ws.Cells.Item(14, 3).Value = 0.08345010231857941
$ws.Cells.Item(14, 4).Value = 0.04993725498158597
$ws.Cells.Item(14, 6).Value = 0.441029643532211
$ws.Cells.Item(14, 7).Value = 0.2807203114955215
$ws.Cells.Item(14, 8).Value = 0.4594234486938475
$ws.Cells.Item(14, 11).Value = 0.4503943955154455
$ws.Cells.Item(14, 14).Value = 1.069530155904538
$ws.Cells.Item(14, 15).Value = 1.404127774164493
$ws.Cells.Item(15, 2).Value = 0.4524109193065726
$ws.Cells.Item(15, 3).Value = 0.08302046241007588
$ws.Cells.Item(15, 4).Value = 0.04944390682578614
$ws.Cells.Item(15, 6).Value = 0.4409863947911816
$ws.Cells.Item(15, 7).Value = 0.2808797298094845
$ws.Cells.Item(15, 8).Value = 0.4598109775682815
$ws.Cells.Item(15, 11).Value = 0.4456443804658932
$ws.Cells.Item(15, 14).Value = 1.070234559177941
$ws.Cells.Item(15, 15).Value = 1.40525679432524
$ws.Cells.Item(16, 2).Value = 0.4252944685719058
$ws.Cells.Item(16, 3).Value = 0.08055545462759994
$ws.Cells.Item(16, 4).Value = 0.0466123889383141
$ws.Cells.Item(16, 6).Value = 0.4408272914535942
$ws.Cells.Item(16, 7).Value = 0.2818619119556018
$ws.Cells.Item(16, 8).Value = 0.4620919848213205
$ws.Cells.Item(16, 11).Value = 0.4183915264318046
$ws.Cells.Item(16, 14).Value = 1.074384649703823
$ws.Cells.Item(16, 15).Value = 1.411996393509654
$ws.Cells.Item(17, 2).Value = 0.4086436589058167
$ws.Cells.Item(17, 3).Value = 0.07904071643493182
$ws.Cells.Item(17, 4).Value = 0.04487152462800736
$ws.Cells.Item(17, 6).Value = 0.4408088835910533
$ws.Cells.Item(17, 7).Value = 0.2825256336797182
$ws.Cells.Item(17, 8).Value = 0.4635450798484371
$ws.Cells.Item(17, 11).Value = 0.4016444216174477
$ws.Cells.Item(17, 14).Value = 1.07703188763989
$ws.Cells.Item(17, 15).Value = 1.41637146829126
$ws.Cells.Item(18, 2).Value = 0.399060996930757
$ws.Cells.Item(18, 3).Value = 0.0781685731906947
$ws.Cells.Item(18, 4).Value = 0.04386884235701416
$ws.Cells.Item(18, 6).Value = 0.4408274060945061
$ws.Cells.Item(18, 7).Value = 0.2829298555797592
$ws.Cells.Item(18, 8).Value = 0.4644006307720616
$ws.Cells.Item(18, 11).Value = 0.3920017384452024
$ws.Cells.Item(18, 14).Value = 1.078591760103549
$ws.Cells.Item(18, 15).Value = 1.41897627376926
$ws.Cells.Item(19, 2).Value = 0.3958155396694565
$ws.Cells.Item(19, 3).Value = 0.0778731273500739
$ws.Cells.Item(19, 4).Value = 0.04352911576101803
$ws.Cells.Item(19, 6).Value = 0.4408386777237467
$ws.Cells.Item(19, 7).Value = 0.2830705750854037
$ws.Cells.Item(19, 8).Value = 0.4646937017776622
$ws.Cells.Item(19, 11).Value = 0.3887351631675529
$ws.Cells.Item(19, 14).Value = 1.079126306614455
$ws.Cells.Item(19, 15).Value = 1.419873393000401
$ws.Cells.Item(20, 2).Value = 0.4104167462923556
$ws.Cells.Item(20, 3).Value = 0.07920205724539642
$ws.Cells.Item(20, 4).Value = 0.04505698622239152
$ws.Cells.Item(20, 6).Value = 0.4408078304659924
$ws.Cells.Item(20, 7).Value = 0.2824526537294503
$ws.Cells.Item(20, 8).Value = 0.4633883495658466
$ws.Cells.Item(20, 11).Value = 0.4034282395057289
$ws.Cells.Item(20, 14).Value = 1.076746229889388
$ws.Cells.Item(20, 15).Value = 1.415896586785166
$ws.Cells.Item(21, 2).Value = 0.4594067986747064
$ws.Cells.Item(21, 3).Value = 0.08365607384274654
$ws.Cells.Item(21, 4).Value = 0.05017375134093527
$ws.Cells.Item(21, 6).Value = 0.4410519519720353
$ws.Cells.Item(21, 7).Value = 0.2806450805700962
$ws.Cells.Item(21, 8).Value = 0.4592386832968032
$ws.Cells.Item(21, 11).Value = 0.4526715747957439
$ws.Cells.Item(21, 14).Value = 1.069194381126529
$ws.Cells.Item(21, 15).Value = 1.403591148854019
$ws.Cells.Item(22, 2).Value = 0.4913702672305078
$ws.Cells.Item(22, 3).Value = 0.0865584126462835
$ws.Cells.Item(22, 4).Value = 0.05350509808953063
$ws.Cells.Item(22, 6).Value = 0.4414706001921829
$ws.Cells.Item(22, 7).Value = 0.2796641809548319
$ws.Cells.Item(22, 8).Value = 0.4567024986215529
$ws.Cells.Item(22, 11).Value = 0.4847592728277448
$ws.Cells.Item(22, 14).Value = 1.064589964085968
$ws.Cells.Item(22, 15).Value = 1.396336310062111
$ws.Cells.Item(23, 2).Value = 0.4743159730454636
$ws.Cells.Item(23, 3).Value = 0.08501019506951479
$ws.Cells.Item(23, 4).Value = 0.05172828397358842
$ws.Cells.Item(23, 6).Value = 0.4412234581793015
$ws.Cells.Item(23, 7).Value = 0.2801693452799938
$ws.Cells.Item(23, 8).Value = 0.4580400253912984
$ws.Cells.Item(23, 11).Value = 0.467642456038476
$ws.Cells.Item(23, 14).Value = 1.067017154599796
$ws.Cells.Item(23, 15).Value = 1.400136277134266
$ws.Cells.Item(24, 2).Value = 0.4096151642352481
$ws.Cells.Item(24, 3).Value = 0.07912911910884191
$ws.Cells.Item(24, 4).Value = 0.04497314474994596
$ws.Cells.Item(24, 6).Value = 0.4408082159227718
$ws.Cells.Item(24, 7).Value = 0.2824855774500818
$ws.Cells.Item(24, 8).Value = 0.4634591445683824
$ws.Cells.Item(24, 11).Value = 0.4026218207032457
$ws.Cells.Item(24, 14).Value = 1.076875257572482
$ws.Cells.Item(24, 15).Value = 1.416111001884815
$ws.Cells.Item(25, 2).Value = 0.3396150646206308
$ws.Cells.Item(25, 3).Value = 0.07275109822649029
$ws.Cells.Item(25, 4).Value = 0.03763398417657982
$ws.Cells.Item(25, 6).Value = 0.4414732901502703
$ws.Cells.Item(25, 7).Value = 0.2858427081848234
$ws.Cells.Item(25, 8).Value = 0.4700604684290752
$ws.Cells.Item(25, 11).Value = 0.3320995651804708
$ws.Cells.Item(25, 14).Value = 1.088933373163954
$ws.Cells.Item(25, 15).Value = 1.436720899642893
